$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 209; this shifts rows 209..278 down to 210..279
$ws.Rows.Item(209).Insert()

# Populate the newly inserted row 209 with a new price-list entry
$ws.Cells.Item(209, 1).Value = 5
$ws.Cells.Item(209, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(209, 3).Value = "Maule"
$ws.Cells.Item(209, 4).Value = 44795
$ws.Cells.Item(209, 5).Value = 7
$ws.Cells.Item(209, 6).Value = "Fruta"
$ws.Cells.Item(209, 7).Value = 100108
$ws.Cells.Item(209, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(209, 9).Value = 100108005
$ws.Cells.Item(209, 10).Value = "Piña"
$ws.Cells.Item(209, 11).Value = "Caramelo"
$ws.Cells.Item(209, 12).Value = "Segunda"
$ws.Cells.Item(209, 13).Value = 540
$ws.Cells.Item(209, 14).Value = 18000
$ws.Cells.Item(209, 15).Value = 18000
$ws.Cells.Item(209, 16).Value = 18000
$ws.Cells.Item(209, 17).Value = "$/caja 14 unidades"
$ws.Cells.Item(209, 18).Value = "Ecuador"
$ws.Cells.Item(209, 19).Value = 1286
$ws.Cells.Item(209, 20).Value = 14

# Make sure the date cell keeps the same date-formatted style as the rest of column D
$ws.Cells.Item(209, 4).NumberFormat = $ws.Cells.Item(210, 4).NumberFormat
